$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.508.05"
$ws.Range("E2").Value = "  +0.22%  "

$ws.Range("D3").Value = "1.915.18"
$ws.Range("E3").Value = "  -0.01%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9995"
$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.39"
$ws.Range("E5").Value = "  +0.61%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9993"
$ws.Range("E6").Value = "  -0.09%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4807"
$ws.Range("E7").Value = "  +2.26%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2882"
$ws.Range("E8").Value = "  +0.60%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06730"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "110.85"
$ws.Range("E10").Value = "  +1.23%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.14"
$ws.Range("E11").Value = "  +3.83%  "

$ws.Range("D12").Value = "1.914.08"
$ws.Range("E12").Value = "  +1.36%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07556"
$ws.Range("E13").Value = "  -2.39%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.250"
$ws.Range("E14").Value = "  -0.78%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6687"
$ws.Range("E15").Value = "  +1.57%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "293.17"
$ws.Range("E16").Value = "  -0.35%  "

$ws.Range("D17").Value = "30.493.53"
$ws.Range("E17").Value = "  +0.22%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9998"
$ws.Range("E18").Value = "  -0.06%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007576"
$ws.Range("E19").Value = "  -0.59%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.92"
$ws.Range("E20").Value = "  -0.11%  "

$ws.Range("D21").Value = "2.161.17"
$ws.Range("E21").Value = "  +1.32%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.477"
$ws.Range("E22").Value = "  +4.35%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9989"
$ws.Range("E23").Value = "  -0.15%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.406"
$ws.Range("E24").Value = "  +3.22%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.472"
$ws.Range("E25").Value = "  +1.46%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.27"
$ws.Range("E26").Value = "  -2.61%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.48"
$ws.Range("E27").Value = "  -5.96%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.122"
$ws.Range("E28").Value = "  +1.86%  "

$ws.Range("E29").Value = "  -0.23%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.401"
$ws.Range("E30").Value = "  +2.56%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.167"
$ws.Range("E31").Value = "  -0.19%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.027"
$ws.Range("E32").Value = "  +1.15%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04970"
$ws.Range("E33").Value = "  -1.46%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7300"
$ws.Range("E34").Value = "  -1.15%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.135"
$ws.Range("E35").Value = "  -1.72%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02056"
$ws.Range("E36").Value = "  -0.43%  "

$ws.Range("B37").Value = "Frax"
$ws.Range("C37").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9995"
$ws.Range("E37").Value = "  +0.04%  "

$ws.Range("B38").Value = "HuobiToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.738"
$ws.Range("E38").Value = "  -0.22%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.670"
$ws.Range("E39").Value = "  -0.31%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "110.94"
$ws.Range("E40").Value = "  +1.48%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.014"
$ws.Range("E41").Value = "  -2.26%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4421"
$ws.Range("E42").Value = "  +3.60%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8638"
$ws.Range("E43").Value = "  -1.01%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.878"
$ws.Range("E44").Value = "  +0.85%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9993"
$ws.Range("E45").Value = "  -0.07%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "68.15"
$ws.Range("E46").Value = "  +1.14%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.315"
$ws.Range("E47").Value = "  +1.61%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "49.06"
$ws.Range("E48").Value = "  -5.01%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.297"
$ws.Range("E49").Value = "  +0.40%  "

$ws.Range("E50").Value = "  +1.55%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.2541"
$ws.Range("E51").Value = "  +3.91%  "
